$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A values (rows 1-33)
$ws.Range("A1").Value = 0.11720183821370256
$ws.Range("A2").Value = -0.0059999999854802866
$ws.Range("A3").Value = -0.0039999999876343395
$ws.Range("A4").Value = -0.0079999999769881924
$ws.Range("A5").Value = -0.0029999999876180183
$ws.Range("A6").Value = -0.0019999999869995122
$ws.Range("A7").Value = -0.0099999999677322648
$ws.Range("A8").Value = -0.0099999999670536965
$ws.Range("A9").Value = -0.0019999999852942096
$ws.Range("A10").Value = -0.001999999984743539
$ws.Range("A11").Value = -0.0029999999823688839
$ws.Range("A12").Value = -0.0034999999810909621
$ws.Range("A13").Value = -0.0034999999808773552
$ws.Range("A14").Value = -0.0079999999702495828
$ws.Range("A15").Value = -0.00099999998699917825
$ws.Range("A16").Value = -0.0019999999848900885
$ws.Range("A17").Value = 0.015767065741737696
$ws.Range("A18").Value = -0.0039999999807287523
$ws.Range("A19").Value = 0.03385085310039404
$ws.Range("A20").Value = -0.0039999999895528049
$ws.Range("A21").Value = -0.0039999999894266836
$ws.Range("A22").Value = -0.0039999999893360894
$ws.Range("A23").Value = -0.0049999999845171672
$ws.Range("A24").Value = -0.019999999947589053
$ws.Range("A25").Value = -0.018188749791983483
$ws.Range("A26").Value = -0.0024999999855559452
$ws.Range("A27").Value = -0.0024999999852788335
$ws.Range("A28").Value = -0.0019999999849558137
$ws.Range("A29").Value = 0.023519860683967053
$ws.Range("A30").Value = -0.059999999845826935
$ws.Range("A31").Value = -0.006999999971444737
$ws.Range("A32").Value = -0.0099999999645294935
$ws.Range("A33").Value = 0.025990885742677605

# Narrow column A width (16.42578125 -> 15.42578125 in OOXML units)
$ws.Columns.Item(1).ColumnWidth = 14.666666666666666
